$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rng = $ws.Range("A2:A49")

# Keep the values stored as plain text (not auto-converted to a date
# serial number) by formatting the cells as Text before assigning ...
$rng.NumberFormat = "@"

# Correct the libraryDate values in column A (rows 2-49): both "2.10.17"
# and "2.1.17" were missing the leading zero; unify them to "02.10.17".
$rng.Value = "02.10.17"

# ... then drop back to the default "Normal" style so the cells keep
# looking exactly as they did before (no explicit style applied to them).
$rng.Style = "Normal"

# Reflect the selection left behind after making the edit.
$ws.Range("A3:A49").Select()
